# Scheduled runner update: refresh market-price / profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on a handful of leve rows
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 736
$ws.Range("I92").Value = 738.1429000000001
$ws.Range("J92").Value = 726
$ws.Range("K92").Value = 738.1429000000001
$ws.Range("L92").Value = 726
$ws.Range("M92").Value = 509.8570999999999
$ws.Range("N92").Value = -3222
$ws.Range("H95").Value = 28500
$ws.Range("J95").Value = 28500
$ws.Range("L95").Value = 28500
$ws.Range("N95").Value = -33992
$ws.Range("H96").Value = 1770.091
$ws.Range("I96").Value = 747
$ws.Range("J96").Value = 2997.8
$ws.Range("K96").Value = 2241
$ws.Range("L96").Value = 8993.400000000001
$ws.Range("M96").Value = -868
$ws.Range("N96").Value = -11739.4
$ws.Range("H97").Value = 125927.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 125927.5
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 377782.5
$ws.Range("N97").Value = -378774.5
$ws.Range("H99").Value = 50250
$ws.Range("I99").Value = 50250
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 150750
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -149252
$ws.Range("H100").Value = 946.875
$ws.Range("I100").Value = 815
$ws.Range("J100").Value = 1166.6666
$ws.Range("K100").Value = 815
$ws.Range("L100").Value = 1166.6666
$ws.Range("M100").Value = -274
$ws.Range("N100").Value = -2248.6666
$ws.Range("H101").Value = 935.75
$ws.Range("I101").Value = 319.33334
$ws.Range("J101").Value = 2785
$ws.Range("K101").Value = 958.0000200000001
$ws.Range("L101").Value = 8355
$ws.Range("M101").Value = 663.9999799999999
$ws.Range("N101").Value = -11599
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
$ws.Range("H106").Value = 2753.2
$ws.Range("I106").Value = 2705.7144
$ws.Range("J106").Value = 2864
$ws.Range("K106").Value = 2705.7144
$ws.Range("L106").Value = 2864
$ws.Range("M106").Value = -2074.7144
$ws.Range("N106").Value = -4126
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27955.125
$ws.Range("I32").Value = 7453.864
$ws.Range("J32").Value = 253469
$ws.Range("K32").Value = 7453.864
$ws.Range("L32").Value = 253469
$ws.Range("M32").Value = -7166.864
$ws.Range("N32").Value = -254043
$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16512
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2190.8076
$ws.Range("I80").Value = 837.4
$ws.Range("J80").Value = 2513.0476
$ws.Range("K80").Value = 837.4
$ws.Range("L80").Value = 2513.0476
$ws.Range("M80").Value = 160.6
$ws.Range("N80").Value = -4509.0476
$ws.Range("H83").Value = 2190.8076
$ws.Range("I83").Value = 837.4
$ws.Range("J83").Value = 2513.0476
$ws.Range("K83").Value = 4187
$ws.Range("L83").Value = 12565.238
$ws.Range("M83").Value = 805
$ws.Range("N83").Value = -22549.238
$ws.Range("H86").Value = 75212.8
$ws.Range("I86").Value = 101755.63
$ws.Range("J86").Value = 2220
$ws.Range("K86").Value = 101755.63
$ws.Range("L86").Value = 2220
$ws.Range("M86").Value = -100632.63
$ws.Range("N86").Value = -4466
$ws.Range("H89").Value = 75212.8
$ws.Range("I89").Value = 101755.63
$ws.Range("J89").Value = 2220
$ws.Range("K89").Value = 508778.15
$ws.Range("L89").Value = 11100
$ws.Range("M89").Value = -503162.15
$ws.Range("N89").Value = -22332
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1214.4667
$ws.Range("I94").Value = 1003
$ws.Range("J94").Value = 1291.3636
$ws.Range("K94").Value = 1003
$ws.Range("L94").Value = 1291.3636
$ws.Range("M94").Value = -552
$ws.Range("N94").Value = -2193.3636
$ws.Range("H134").Value = 1422
$ws.Range("I134").Value = 1422
$ws.Range("J134").Value = 1422
$ws.Range("K134").Value = 4266
$ws.Range("L134").Value = 4266
$ws.Range("M134").Value = -1731
$ws.Range("N134").Value = -9336
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1100
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -5717
$ws.Range("N32").Value = -2066
$ws.Range("H34").Value = 705
$ws.Range("I34").Value = 100.5
$ws.Range("J34").Value = 877.7143
$ws.Range("K34").Value = 301.5
$ws.Range("L34").Value = 2633.1429
$ws.Range("M34").Value = -217.5
$ws.Range("N34").Value = -2801.1429
$ws.Range("H107").Value = 604985.5
$ws.Range("I107").Value = 946
$ws.Range("K107").Value = 2838
$ws.Range("M107").Value = -918
$ws.Range("H131").Value = 777.89
$ws.Range("J131").Value = 790.51044
$ws.Range("L131").Value = 2371.53132
$ws.Range("N131").Value = -12451.53132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 6867.4116
$ws.Range("I45").Value = 6666.6665
$ws.Range("J45").Value = 6910.4287
$ws.Range("K45").Value = 6666.6665
$ws.Range("L45").Value = 6910.4287
$ws.Range("M45").Value = -6259.6665
$ws.Range("N45").Value = -7724.4287
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2163.6956
$ws.Range("I122").Value = 1230.5
$ws.Range("J122").Value = 2661.4
$ws.Range("K122").Value = 3691.5
$ws.Range("L122").Value = 7984.200000000001
$ws.Range("M122").Value = -1241.5
$ws.Range("N122").Value = -12884.2
$ws.Range("H126").Value = 1795.9
$ws.Range("I126").Value = 1556
$ws.Range("J126").Value = 2035.8
$ws.Range("K126").Value = 4668
$ws.Range("L126").Value = 6107.4
$ws.Range("M126").Value = -2198
$ws.Range("N126").Value = -11047.4
$ws.Range("H132").Value = 2421.239
$ws.Range("I132").Value = 2748.0667
$ws.Range("J132").Value = 1808.4375
$ws.Range("K132").Value = 8244.2001
$ws.Range("L132").Value = 5425.3125
$ws.Range("M132").Value = -5714.2001
$ws.Range("N132").Value = -10485.3125
